# Update generated counts in "展览" (Exhibition) and "全部类型" (All Types) sheets:
#   F2: 533 -> 532
#   F4: 25  -> 26

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 532
    $ws.Range("F4").Value = 26
}
